$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook rows were reordered/edited per the commit; update each
# affected cell to its new value. A leading escaped apostrophe forces
# Excel to store the value as literal text (avoids 'True' -> Boolean).
$ws.Range("A2").Value = '''Observation'
$ws.Range("B2").Value = '''RespiratoryRate'
$ws.Range("C2").Value = '''Første respirationsfrekvens måling.'
$ws.Range("D4").Value = '''ResultValue'
$ws.Range("E4").Value = '''Integer'
$ws.Range("G4").Value = '''Greater than or equal to: 0'
$ws.Range("B5").Value = '''PulseRate'
$ws.Range("C5").Value = '''Første pulsmåling.'
$ws.Range("E7").Value = '''Integer'
$ws.Range("G7").Value = '''Greater than or equal to: 0'
$ws.Range("B8").Value = '''PainEvaluation'
$ws.Range("C8").Value = '''Første smertemåling.'
$ws.Range("D9").Value = '''ResultCode'
$ws.Range("E9").Value = '''str, Enum'
$ws.Range("F9").Value = '''Se webservice dokumentation.'
$ws.Range("G9").Value = '''Enums/Udfald: | "vas" | "nrs" | "vrs" | "andet" | '
$ws.Range("D10").Value = '''DateTime'
$ws.Range("E10").Value = '''str'
$ws.Range("G10").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("D11").Value = '''ResultValue'
$ws.Range("E11").Value = '''Float'
$ws.Range("F11").Value = '''Værdi for smertescore'
$ws.Range("G11").Value = '''Greater than or equal to: 0'
$ws.Range("H11").Value = '''True'
$ws.Range("D12").Value = '''Note'
$ws.Range("E12").Value = '''String'
$ws.Range("F12").Value = '''Bruges ved ''andet'' type af smertescore'
$ws.Range("A13").Value = '''Observation'
$ws.Range("B13").Value = '''SystolicBloodPressure'
$ws.Range("C13").Value = '''Første systoliske blodtryksmåling'
$ws.Range("D14").Value = '''DateTime'
$ws.Range("E14").Value = '''str'
$ws.Range("G14").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("H14").Value = '''True'
$ws.Range("D15").Value = '''ResultValue'
$ws.Range("E15").Value = '''Integer'
$ws.Range("G15").Value = '''Greater than or equal to: 0'
$ws.Range("A16").Value = '''Observation'
$ws.Range("B16").Value = '''OxygenDemand'
$ws.Range("C16").Value = '''Første iltbehov.'
$ws.Range("D17").Value = '''ResultCode'
$ws.Range("E17").Value = '''str, Enum'
$ws.Range("G17").Value = '''Enums/Udfald: | "y" | "n" | '
$ws.Range("D18").Value = '''DateTime'
$ws.Range("E18").Value = '''str'
$ws.Range("G18").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("H18").Value = '''True'
$ws.Range("B19").Value = '''AVPUScale'
$ws.Range("C19").Value = '''Første AVPU skala'
$ws.Range("D20").Value = '''ResultCode'
$ws.Range("E20").Value = '''str, Enum'
$ws.Range("G20").Value = '''Enums/Udfald: | "a" | "v" | "p" | "u" | '
$ws.Range("D21").Value = '''DateTime'
$ws.Range("E21").Value = '''str'
$ws.Range("G21").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("B22").Value = '''GlasgowComaScale'
$ws.Range("C22").Value = '''Første Glasgow Coma Scale måling.'
$ws.Range("D23").Value = '''DateTime'
$ws.Range("E23").Value = '''str'
$ws.Range("G23").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("D24").Value = '''ResultValue'
$ws.Range("E24").Value = '''Integer'
$ws.Range("G24").Value = '''Greater than or equal to: 3 | Less than or equal to: 15'
$ws.Range("B25").Value = '''Triage'
$ws.Range("C25").Value = '''Første Triage kategorisering.'
$ws.Range("G26").Value = '''Enums/Udfald: | "r" | "o" | "y" | "g" | "b" | '
$ws.Range("B31").Value = '''OxygenSaturation'
$ws.Range("C31").Value = '''Første iltmætningsmåling.'
$ws.Range("E33").Value = '''Float'
$ws.Range("F33").Value = '''Enhed = %'
$ws.Range("G33").Value = '''Greater than or equal to: 0 | Less than or equal to: 100'
$ws.Range("A34").Value = '''Contact'
$ws.Range("B34").Value = '''NoteType'
$ws.Range("C34").Value = '''Første afsluttet lægenotat'
$ws.Range("D35").Value = '''DateTime'
$ws.Range("E35").Value = '''str'
$ws.Range("F35").Value = '''Dato for første afsluttet lægenotat'
$ws.Range("G35").Value = '''Date format: | DD-MM-YYYY HH:MM:SS |'
$ws.Range("D36").Value = '''Note'
$ws.Range("E36").Value = '''String'
$ws.Range("F36").Value = '''Notattype for første afsluttet lægenotat'

# Clear cells that no longer hold data at these positions
$ws.Range("F3").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("F7").Value = $null
$ws.Range("A11").Value = $null
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("G12").Value = $null
$ws.Range("H12").Value = $null
$ws.Range("D13").Value = $null
$ws.Range("E13").Value = $null
$ws.Range("G13").Value = $null
$ws.Range("H13").Value = $null
$ws.Range("A14").Value = $null
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = $null
$ws.Range("F15").Value = $null
$ws.Range("D16").Value = $null
$ws.Range("E16").Value = $null
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("F17").Value = $null
$ws.Range("F18").Value = $null
$ws.Range("G36").Value = $null
